$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.533.37"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "1.953.56"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'243.98"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'0.621"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D7").Value = "'58.47"
$ws.Range("E7").Value = "  -2.94%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").Value = "'55.82"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +5.46%  "
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'21.97"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").Value = "'0.828"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("D15").Value = "2.233.45"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "'13.68"
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("D17").Value = "'5.25"
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "1.937.00"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "36.428.36"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").Value = "'69.94"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").Value = "0.0₃0864"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "'229.88"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("E23").Value = "  -3.06%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  -3.70%  "
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("E27").Value = "  -5.89%  "
$ws.Range("D28").Value = "'162.16"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").Value = "'0.138"
$ws.Range("E29").Value = "  +8.55%  "
$ws.Range("D30").Value = "'19.54"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "'1.17"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("D34").Value = "'0.0633"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("D35").Value = "'4.30"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("D36").Value = "'6.33"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "'2.16"
$ws.Range("E39").Value = "  -5.86%  "
$ws.Range("D40").Value = "'3.05"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("D41").Value = "'0.0989"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("E43").Value = "  -4.00%  "
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("D45").Value = "'16.08"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "1.369.48"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("E47").Value = "  -4.70%  "
$ws.Range("D48").Value = "'88.45"
$ws.Range("E48").Value = "  -4.43%  "
$ws.Range("D49").Value = "'7.17"
$ws.Range("E49").Value = "  -4.45%  "
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "'45.91"
$ws.Range("E51").Value = "  +3.93%  "
